# "6 hours by turn fix" - shift the afternoon schedule rows down by one slot
# (each turn becomes 50 minutes later starting at the lunch break) and add
# three new time slots (16:40, 17:30, 18:20) at the end of the day so the
# teacher's day covers a full 6 hours per turn.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Morning block: a couple of classes move between columns/slots --
$ws.Range("E3").Value  = "-"

$ws.Range("D4").Value  = "-"
$ws.Range("E4").Value  = "MCT-1A-Tecnologia dos Materiais"
$ws.Range("F4").Value  = "-"

$ws.Range("C6").Value  = "MEC-1A-Tecnologia dos Materiais"
$ws.Range("F6").Value  = "-"

$ws.Range("D7").Value  = "MCT-1A-Tecnologia dos Materiais"
$ws.Range("E7").Value  = "MEC-1A-Tecnologia dos Materiais"

# -- Row 8 ("Almoço") no longer holds the lunch break; it becomes a normal slot --
$ws.Range("B8").Value  = "-"
$ws.Range("C8").Value  = "-"
$ws.Range("D8").Value  = "-"
$ws.Range("E8").Value  = "-"
$ws.Range("F8").Value  = "-"

# -- Row 9 becomes the lunch break, starting earlier at 12:20 --
$ws.Range("A9").Value  = "12:20"
$ws.Range("B9").Value  = "Almoço"
$ws.Range("C9").Value  = "Almoço"
$ws.Range("D9").Value  = "Almoço"
$ws.Range("E9").Value  = "Almoço"
$ws.Range("F9").Value  = "Almoço"

# -- Rows 10/11 keep "-" content but shift start times earlier by one slot --
$ws.Range("A10").Value = "13:00"
$ws.Range("A11").Value = "13:50"

# -- Row 12 shifts time and stops being "Intervalo" --
$ws.Range("A12").Value = "14:40"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"

# -- Row 13 shifts time and becomes "Intervalo" --
$ws.Range("A13").Value = "15:30"
$ws.Range("B13").Value = "Intervalo"
$ws.Range("C13").Value = "Intervalo"
$ws.Range("D13").Value = "Intervalo"
$ws.Range("E13").Value = "Intervalo"
$ws.Range("F13").Value = "Intervalo"

# -- Row 14 shifts time, keeps "-" content --
$ws.Range("A14").Value = "15:50"

# -- New rows 15-17 appended at the end of the day --
$ws.Range("A15").Value = "16:40"
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"

$ws.Range("A16").Value = "17:30"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"

$ws.Range("A17").Value = "18:20"
